# Applies the "began writing functions" data-file update to
# tutorial/data/transport_model_v1/transport_model_data.xlsx
#
# Summary of the change:
#  - "parameters" sheet: insert a new parameter row ("Energy storage power",
#    offset index 18, value 0, units kW, note "No LIB in a diesel tractor")
#    right after "Energy storage weight" / before "Electric drive size".
#    All rows below shift down one row and their Offset (column D) index
#    is bumped by one to stay sequential.
#  - "results" sheet: append a new "Energy consumption" metric row (units gge).
#  - "indices" sheet: append a matching new "Energy consumption" metric index row.

$wb = $excel.ActiveWorkbook

$wsParams = $wb.Worksheets.Item("parameters")
$wsResults = $wb.Worksheets.Item("results")
$wsIndices = $wb.Worksheets.Item("indices")

# ---------------------------------------------------------------------
# 1. results sheet — append "Energy consumption" metric row
#    (added to the shared-string table first, so "Energy consumption"
#    sorts ahead of "Energy storage power", matching the saved file)
# ---------------------------------------------------------------------
$wsResults.Range("A7").Value = "Class 8 Diesel Tractor"
$wsResults.Range("B7").Value = "Metric"
$wsResults.Range("C7").Value = "Energy consumption"
$wsResults.Range("D7").Value = "gge"

$wsResults.Columns.Item(3).ColumnWidth = 18.5

# ---------------------------------------------------------------------
# 2. indices sheet — append matching "Energy consumption" index row
# ---------------------------------------------------------------------
$wsIndices.Range("A15").Value = "Class 8 Diesel Tractor"
$wsIndices.Range("B15").Value = "Metric"
$wsIndices.Range("C15").Value = "Energy consumption"
$wsIndices.Range("D15").Value = 4

# ---------------------------------------------------------------------
# 3. parameters sheet — insert new "Energy storage power" row at row 20
# ---------------------------------------------------------------------
$wsParams.Rows.Item(20).Insert()

$wsParams.Range("A20").Value = "Class 8 Diesel Tractor"
$wsParams.Range("B20").Value = "Reference"
$wsParams.Range("C20").Value = "Energy storage power"
$wsParams.Range("D20").Value = 18
$wsParams.Range("E20").Value = 0
$wsParams.Range("F20").Value = "kW"
$wsParams.Range("G20").Value = "No LIB in a diesel tractor"

# Every row that was pushed down (old rows 20-28, now 21-29) keeps its data
# but its sequential Offset index (column D) needs to move up by one.
for ($r = 21; $r -le 29; $r++) {
    $cell = $wsParams.Cells.Item($r, 4)
    $cell.Value = $cell.Value2 + 1
}

# ---------------------------------------------------------------------
# 4. Restore on-screen selections to match the saved view state
# ---------------------------------------------------------------------
$wsIndices.Range("A14:A15").Select()
$wsResults.Range("C14").Select()
$wsParams.Range("C20").Select()

$wsDesigns = $wb.Worksheets.Item("designs")
$wsDesigns.Range("D3:D8").Select()

Write-Output "edit applied"
